# Cuppa: Removed unused slides in figures.pptx
#
# 1. Delete the first slide (sldId 269 - single full-bleed picture) and the
#    last slide (sldId 272 - single picture) leaving the three slides in
#    between (sldIds 270, 271, 268) in their original relative order.
# 2. The footer "date last updated" field on the slide master and every
#    slide layout advances from 06/11/2023 to 19/01/2024 (re-saving the
#    deck re-stamps the auto date placeholder).

$p = $ppt.ActivePresentation

# --- 1. Remove the two picture-only slides (sldId 269 and 272) -----------

$idsToDelete = @(269, 272)

$indexesToDelete = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sld = $p.Slides.Item($i)
    if ($idsToDelete -contains $sld.SlideID) {
        $indexesToDelete.Add($i) | Out-Null
    }
}

# Delete from the highest index down so earlier indexes stay valid.
for ($k = $indexesToDelete.Count - 1; $k -ge 0; $k--) {
    $idx = $indexesToDelete[$k]
    $p.Slides.Item($idx).Delete()
}

# --- 2. Refresh the footer date placeholder text --------------------------

$newDate = "19/01/2024"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
